{"js": "// Design notes for branching:\n// \"... a simple timeline involving two sessions alongside the data ...\"\n// becomes\n// \"... a simple timeline involving two completed sessions alongside the data ...\"\n//\n// Locate the (unique) sentence fragment and insert the word \"completed \"\n// immediately before \"sessions\".\n\nconst body = context.document.body;\nconst results = body.search(\"sessions alongside the data\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase not found\");\n}\n\nconst hit = results.items[0];\nhit.getRange(\"Start\").insertText(\"completed \", \"Before\");\nawait context.sync();\n", "ps1": "# Design notes for branching:\n# \"... a simple timeline involving two sessions alongside the data ...\"\n# becomes\n# \"... a simple timeline involving two completed sessions alongside the data ...\"\n#\n# Find the (unique) sentence fragment and insert the word \"completed \"\n# immediately before it, preserving the surrounding run formatting.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"sessions alongside the data\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\n\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $rng.InsertBefore(\"completed \")\n}\n"}
